$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1 -and $parts[0] -eq "System") {
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + ,"System"
        $cell.Value = [string]::Join(", ", $newParts)
    }
}
